# Insert a new data row at row 330 (pushes existing rows 330-369 down to 331-370)
# and populate it with the new Jengibre price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("330:330").Insert()

$ws.Cells.Item(330, 1).Value = 10
$ws.Cells.Item(330, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(330, 3).Value = "La Araucanía"
$ws.Cells.Item(330, 4).Value = 45218
$ws.Cells.Item(330, 5).Value = 9
$ws.Cells.Item(330, 6).Value = 100114007
$ws.Cells.Item(330, 7).Value = "Jengibre"
$ws.Cells.Item(330, 8).Value = "Sin especificar"
$ws.Cells.Item(330, 9).Value = "Primera"
$ws.Cells.Item(330, 10).Value = 100
$ws.Cells.Item(330, 11).Value = 26000
$ws.Cells.Item(330, 12).Value = 26000
$ws.Cells.Item(330, 13).Value = 26000
$ws.Cells.Item(330, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(330, 15).Value = "Perú"
$ws.Cells.Item(330, 16).Value = 2000
$ws.Cells.Item(330, 17).Value = 13
$ws.Cells.Item(330, 18).Value = "Hortaliza"
